$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header text in F1
$ws.Range("F1").Value = "ultima comando"

# Add new data row 2
$ws.Range("A2").Value = 1972887851
$ws.Range("B2").Value = "If they had finished their work earlier, they could have gone to the party."
$ws.Range("C2").Value = "Se eles tivessem terminado o trabalho mais cedo, eles poderiam ter ido à festa."
$ws.Range("D2").Value = "Intermediário"
$ws.Range("E2").Value = 86
$ws.Range("F2").Value = "/OK"
$ws.Range("G2").Value = "Frase"

# New row should not inherit the header's bold/centered style
$ws.Range("A2:G2").Style = "Normal"

# Move selection to match final state
$null = $ws.Range("F8").Select()
